# Updated charts with data for Dec. 2025 + plus annotation changes
#
# Applies:
#   1. Metadata text updates in sharedStrings (Series Name list, "To month").
#   2. New "2025-12" period row (row 69) with that month's tourist-arrival data.
#   3. A fresh trailing blank row (row 70) matching the sheet's previous
#      "next empty period" row, so the table keeps one blank row below data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Metadata text cells (row 13 = "Series Name: ...", row 19 = "To month: ..")
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Series Name: Asia,total, Uzbekistan, Azerbaijan, United Arab Emirates, Indonesia, Armenia, Georgia, India, Hong Kong, Philippines, Vietnam, Turkey, Taiwan, Japan, Jordan, Malaysia, Nepal, China, Singapore, South Korea, Kazakhstan, Cambodia, Cyprus, Thailand, Asia,Othercountries, Africa,total, Uganda, Ethiopia, Ghana, SouthAfrica, Zimbabwe, IvoryCoast, Tanzania, Mauritius, Egypt, Morocco, Nigeria, Kenya, Rwanda, Tunisia, Africa,Other countries, Europe,total, Austria, Ukraine, Italy, Ireland, Albania, Estonia, Bulgaria, Belgium, Belarus, Germany, Denmark, Netherlands, Hungary, United Kingdom, Greece, Luxembourg, Latvia, Lithuania, Moldova, Montenegro, Malta, Macedonia, Norway, Slovenia, Slovakia, Spain, Serbia, Poland, Portugal, Finland, Czech Republic, France, Croatia, Romania, Russian Federation, Sweden, Switzerland, Europe,Other countries, NorthAmerica,total, United States, Canada, Central America, total, ElSalvador, Guatemala, Honduras, Dominican Republic, Mexico, Panama, Costa Rica, Centra lAmerica,Othercoun, South America,total, Uruguay, Ecuador, Argentina, Bolivia, Brazil, Paraguay, Peru, Chile, Colombia, South America,Other countr, Oceania, total, Australia, NewZealand, Oceania,Other countries, Unclas sified countries"

$ws.Range("A19").Value = "To month: 12"

# ---------------------------------------------------------------------------
# 2. Push the previously-blank trailing row (69) down to row 70 so a blank
#    "next period" row still trails the data after we fill row 69 in.
# ---------------------------------------------------------------------------
$lastCol = 109   # column DE
$srcRow = 69
$dstRow = 70

for ($c = 1; $c -le $lastCol; $c++) {
    $srcCell = $ws.Cells.Item($srcRow, $c)
    $dstCell = $ws.Cells.Item($dstRow, $c)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 3. Fill row 69 with the new "2025-12" period.
# ---------------------------------------------------------------------------
$ws.Cells.Item(69, 1).Value = "2025-12"

$decValues = @(9.6,0.1,0.2,0,1.8,0.1,0.4,1.2,0,0.5,0,0.3,0.2,0.3,0.6,0.2,0,1,1.2,0.4,0.1,0,0.4,0.1,0.5,3.4,0,0.1,0,2.1,0,0,0.1,0,0.1,0.1,0.4,0.1,0,0,0.2,63.4,1.1,2.5,3.7,0.2,0,0.1,0.4,1.6,0.4,5.1,0.3,1.4,0.7,10,0.8,0,0.2,0.4,0.3,0,0,0.1,0.3,0.1,0.2,2.1,0.2,0.9,0.5,0.2,0.6,17.9,0.1,2.4,6.2,0.6,1.6,0.1,52.1,48.3,3.7,2.6,0,0.1,0.1,0,1.9,0.3,0.1,0.1,5.2,0.1,0.1,1.7,0,2.4,0,0.2,0.3,0.3,0,2.3,2.1,0.1,0,0.3)

# Columns B..DE (2..109) get the numeric data, with the same #,##0.0 number
# style used by every other data row -- copy that formatting from the row
# directly above (row 68) before writing the values in.
$ws.Range("B68:DE68").Copy()
$ws.Range("B69:DE69").PasteSpecial(-4122)

$col = 2
foreach ($v in $decValues) {
    $ws.Cells.Item(69, $col).Value = $v
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 4. Re-stamp the header picture's internal name (new GUID on every re-embed).
# ---------------------------------------------------------------------------
if ($ws.Shapes.Count -ge 1) {
    $ws.Shapes.Item(1).Name = "Picture c7b46d57-52d3-4f5f-a544-9c6d35f454c4"
}

Write-Output "Row 69 (2025-12) populated; row 70 blank row appended."
